$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.888.36'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.646.91'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0621'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.77'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.648.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.33'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.928.63'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.76'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +7.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.39'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.43'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.18'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.92'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.35'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.87'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0512'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.246.41'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.830'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.05%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.808'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.790.60'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.10'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.90'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.57'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.59'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.45%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.61'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.11%  '
